$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '250.60'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.93'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.437'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05626'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.418'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.377'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8152'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9155'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1429'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07501'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03128'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03097'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09347'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.560'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001607'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04752'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0005797'
$ws.Range("E18").Value = '17OneONEWorstin24h'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006419'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.004995'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.001033'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0001501'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.708'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.190'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04019'

$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006785'
$ws.Range("E41").Value = '40KickTokenKICK'

$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1068'
$ws.Range("E42").Value = '41BKEXTokenBKK'

$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002754'
$ws.Range("E43").Value = '42CEJICEJI'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007877'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.2376'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002102'
